$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Upload" column (G) ---
$ws.Range("G1").Value = "Upload"
$ws.Range("G2").Value = 1

# G1 gets the same header look as the other header cells (font/fill from
# F1's "Uploading type" header), but with a left+right only thin border
# instead of a full box.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Borders.LineStyle = 0
$ws.Range("G1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("G1").Borders.Item(10).LineStyle = 1  # xlEdgeRight

# PasteSpecial only touches formatting, but make sure the header text is
# still exactly "Upload" (not whatever was copied from F1).
$ws.Range("G1").Value = "Upload"

# Move the active cell/selection from E2 to A2.
$ws.Range("A2").Select()
